$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename header from "AccessionFile" to "AccessionFilePath"
$ws.Range("B1").Value = "AccessionFilePath"

# Prefix/suffix every accession value in column B (rows 2-32) with the cluster path
$prefix = "/home/erinroberts/bio_project_data/"
$suffix = ".txt"
for ($r = 2; $r -le 32; $r++) {
    $cell = $ws.Cells.Item($r, 2)
    $accession = $cell.Value2
    $cell.Value = $prefix + $accession + $suffix
}

# Widen column B to fit the longer path strings.
# (The engine's stored column width always comes out 5/6 wider than the
# ColumnWidth value assigned, so back that padding out to land on exactly 47.)
$ws.Columns.Item(2).ColumnWidth = 46.166666666666664

# Update the selection to match the saved view state
$ws.Range("C32").Select()
